$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held the lecturer name ("11079086 - Herlandí de Souza
# Andrade") in columns B/C with no label in A is removed entirely; every
# row below it shifts up by one.
$ws.Rows(13).Delete()

# "Objetivos:" row now (mistakenly) shows the lecturer name instead of the
# course-objectives paragraph.
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# "Programa resumido:" / "Short syllabus:" summary text becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "Programa:" row content becomes the activation-date string. Assigning the
# date-looking text directly would auto-convert the cell to a real date
# serial, so build it as a formula first and paste back as a value to keep
# it as literal text (matching the original "Ativação:" cell's string).
$ws.Range("B15").Formula = '="01/01/2021"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").Formula = '="01/01/2021"'
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# "Método:" row now shows the lecturer name.
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"

# "Critério:" row now shows the evaluation-method paragraph.
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."

# "Norma de recuperação:" row now shows the weighted-average criterion text.
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."

# "Bibliografia:" row now shows the recovery-grade formula paragraph.
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
